# Updated the ESR locators
# Replace the old locator email (wakram@dacgroup.com) with the new one
# (adevaraj@dacgroup.com) for every scenario row, and make sure the
# Scenario4/No Email row (C4) also gets a "mailto" hyperlink + the
# Hyperlink cell style, matching the other rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newEmail = "adevaraj@dacgroup.com"
$newMailto = "mailto:" + $newEmail

# Remove all existing hyperlinks on the sheet so they can be re-created
# against the new address (the object model always reassigns hyperlink
# relationship ids, so drop-and-recreate keeps things consistent).
$ws.Hyperlinks.Delete()

# Update the displayed text for every locator cell in column C.
$ws.Range("C2").Value = $newEmail
$ws.Range("C3").Value = $newEmail
$ws.Range("C4").Value = $newEmail
$ws.Range("C5").Value = $newEmail

# Re-create the hyperlinks in the same order as before (C2, C5, C3) and
# add the new one for C4.
$ws.Hyperlinks.Add($ws.Range("C2"), $newMailto)
$ws.Hyperlinks.Add($ws.Range("C5"), $newMailto)
$ws.Hyperlinks.Add($ws.Range("C3"), $newMailto)
$ws.Hyperlinks.Add($ws.Range("C4"), $newMailto)

# Make sure all four cells use the built-in Hyperlink style (C4 didn't
# have it before).
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C3").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"

# Match the final selection left behind in the workbook.
$ws.Range("F6").Select()
